# Converted Employer Integration Tests to XUnit
# Employer-MissingMandatory.xlsx: blank out the (Do Not Modify) Account id
# on row 2 (so the row has no pre-existing CRM account id) and stamp the
# Company Name cell with the scenario name, matching the new file name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# (Do Not Modify) Account - remove the existing GUID entirely so the row
# represents a "new" record with no linked account.
$ws.Range("A2").Clear()

# Company Name - identify the scenario directly in the data row.
$ws.Range("D2").Value = "Employer-MissingMandatory"

# Widen Company Name column so the new text isn't clipped, and re-select
# the now-relevant cell as the active cell on the sheet.
$ws.Columns.Item(4).ColumnWidth = 25.86
[void]$ws.Range("D3").Select()
